# Daily attendance processing - sort the "Recorded By" list (column G)
# alphabetically (ordinal / case-sensitive order, uppercase before
# lowercase) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = "Recorded By"
    $val = $cell.Value2

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = $val -split ','

    $list = New-Object System.Collections.Generic.List[string]
    foreach ($p in $parts) {
        [void]$list.Add($p.Trim())
    }

    $list.Sort([System.StringComparer]::Ordinal)

    $newVal = [string]::Join(', ', $list)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
